$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab
$ws.Name = "UniformA"

# Add new row 16, mirroring row 15's shape (A=index, B=label, C:P=1)
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"

for ($col = 3; $col -le 16; $col++) {
    $ws.Cells.Item(16, $col).Value = 1
}

# Copy the formatting (bold/border/centered style) from A15 onto A16
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
